$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values stored as text (e.g. "23.832.76", "1.001").
# Force text format before assignment so Excel does not reinterpret these
# strings as numbers (which would drop formatting like trailing zeros or
# collapse multi-dot "thousand.thousand.decimal" strings).

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "23.832.76"
$ws.Cells.Item(2, 5).Value = "  -2.26%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.620.51"
$ws.Cells.Item(3, 5).Value = "  -1.96%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.57%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "308.55"
$ws.Cells.Item(5, 5).Value = "  -1.18%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  -0.48%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3941"
$ws.Cells.Item(7, 5).Value = "  +0.65%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3850"
$ws.Cells.Item(8, 5).Value = "  -1.53%  "

$ws.Cells.Item(9, 5).Value = "  -0.44%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "49.60"
$ws.Cells.Item(10, 5).Value = "  -2.02%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.354"
$ws.Cells.Item(11, 5).Value = "  -1.93%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08481"
$ws.Cells.Item(12, 5).Value = "  -0.65%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "23.77"
$ws.Cells.Item(13, 5).Value = "  -4.85%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.068"
$ws.Cells.Item(14, 5).Value = "  -1.66%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.612"
$ws.Cells.Item(15, 5).Value = "  +0.05%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.00001285"
$ws.Cells.Item(16, 5).Value = "  -1.32%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.619.62"
$ws.Cells.Item(17, 5).Value = "  -2.39%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "93.93"
$ws.Cells.Item(18, 5).Value = "  +0.88%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06940"
$ws.Cells.Item(19, 5).Value = "  -0.26%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "20.03"
$ws.Cells.Item(20, 5).Value = "  -4.76%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.829"
$ws.Cells.Item(21, 5).Value = "  -2.38%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.000"
$ws.Cells.Item(22, 5).Value = "  -0.53%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "13.43"
$ws.Cells.Item(23, 5).Value = "  -2.61%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "23.817.82"
$ws.Cells.Item(24, 5).Value = "  -2.30%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.493"
$ws.Cells.Item(25, 5).Value = "  +6.27%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.830"
$ws.Cells.Item(26, 5).Value = "  +2.17%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "22.26"
$ws.Cells.Item(27, 5).Value = "  -1.89%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "157.05"
$ws.Cells.Item(28, 5).Value = "  -0.94%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "140.52"
$ws.Cells.Item(29, 5).Value = "  -2.77%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "5.304"
$ws.Cells.Item(30, 5).Value = "  -8.04%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.802"
$ws.Cells.Item(31, 5).Value = "  -3.53%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "2.483"
$ws.Cells.Item(32, 5).Value = "  -1.60%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.792.71"
$ws.Cells.Item(33, 5).Value = "  -2.57%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.08138"
$ws.Cells.Item(34, 5).Value = "  -0.69%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.9884"
$ws.Cells.Item(35, 5).Value = "  -2.90%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "6.634"
$ws.Cells.Item(36, 5).Value = "  -2.95%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.02893"
$ws.Cells.Item(37, 5).Value = "  -4.52%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.2672"
$ws.Cells.Item(38, 5).Value = "  -3.27%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.09167"
$ws.Cells.Item(39, 5).Value = "  -4.18%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "10.40"
$ws.Cells.Item(40, 5).Value = "  +2.33%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "13.71"
$ws.Cells.Item(41, 5).Value = "  +3.29%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.426"
$ws.Cells.Item(42, 5).Value = "  -4.56%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.7524"
$ws.Cells.Item(43, 5).Value = "  -3.17%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "15.92"
$ws.Cells.Item(44, 5).Value = "  -2.27%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.6950"
$ws.Cells.Item(45, 5).Value = "  -0.81%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.469"
$ws.Cells.Item(46, 5).Value = "  -3.21%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "4.072"

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.08248"
$ws.Cells.Item(49, 5).Value = "  -3.57%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "135.64"
$ws.Cells.Item(50, 5).Value = "  -0.84%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.200"
$ws.Cells.Item(51, 5).Value = "  -7.89%  "

